# Updated loading-percent results for the "380 kV" case run (Case_4_174).
# Columns B:E, G, M:O hold the per-line loading percentages for each of the
# 24 scenario rows (F, H:L stay 0 and are untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.0594409910526
$ws.Range("C2").Value = 8.146183392977051
$ws.Range("D2").Value = 5.974586842889962
$ws.Range("E2").Value = 11.16093018045487
$ws.Range("G2").Value = 3.618073060271496
$ws.Range("M2").Value = 14.88817371009976
$ws.Range("N2").Value = 17.06051441188715
$ws.Range("O2").Value = 21.88661657444765
$ws.Range("B3").Value = 12.51247904771413
$ws.Range("C3").Value = 7.687364648947827
$ws.Range("D3").Value = 5.854270157744292
$ws.Range("E3").Value = 11.0445146453906
$ws.Range("G3").Value = 3.620712390535023
$ws.Range("M3").Value = 14.62109935982
$ws.Range("N3").Value = 17.12923465148288
$ws.Range("O3").Value = 21.8559206580594
$ws.Range("B4").Value = 12.1665497709628
$ws.Range("C4").Value = 7.389757060844937
$ws.Range("D4").Value = 5.780982521354287
$ws.Range("E4").Value = 10.97656079946481
$ws.Range("G4").Value = 3.622418115142265
$ws.Range("M4").Value = 14.45850344182482
$ws.Range("N4").Value = 17.17332591930335
$ws.Range("O4").Value = 21.84395831400555
$ws.Range("B5").Value = 12.02326189337851
$ws.Range("C5").Value = 7.264526866488651
$ws.Range("D5").Value = 5.751310998567079
$ws.Range("E5").Value = 10.94978513446981
$ws.Range("G5").Value = 3.62313470106179
$ws.Range("M5").Value = 14.39268571730198
$ws.Range("N5").Value = 17.19177207230505
$ws.Range("O5").Value = 21.84081509449519
$ws.Range("B6").Value = 11.99933591816104
$ws.Range("C6").Value = 7.243495234920383
$ws.Range("D6").Value = 5.746397123751406
$ws.Range("E6").Value = 10.94539520101144
$ws.Range("G6").Value = 3.623254989678489
$ws.Range("M6").Value = 14.38178609440527
$ws.Range("N6").Value = 17.19486399684921
$ws.Range("O6").Value = 21.84039771704964
$ws.Range("B7").Value = 12.16462642074138
$ws.Range("C7").Value = 7.388084101855716
$ws.Range("D7").Value = 5.780581514969801
$ws.Range("E7").Value = 10.97619594706788
$ws.Range("G7").Value = 3.622427692163793
$ws.Range("M7").Value = 14.45761389165195
$ws.Range("N7").Value = 17.17357275056983
$ws.Range("O7").Value = 21.8439089136482
$ws.Range("B8").Value = 12.87306641465725
$ws.Range("C8").Value = 7.99130359534475
$ws.Range("D8").Value = 5.933007260760262
$ws.Range("E8").Value = 11.12007672197195
$ws.Range("G8").Value = 3.618965469832824
$ws.Range("M8").Value = 14.79585340832007
$ws.Range("N8").Value = 17.08381647280099
$ws.Range("O8").Value = 21.87460429401055
$ws.Range("B9").Value = 14.1739214388166
$ws.Range("C9").Value = 9.046724136511267
$ws.Range("D9").Value = 6.234581577821479
$ws.Range("E9").Value = 11.42879913342261
$ws.Range("G9").Value = 3.612848441642258
$ws.Range("M9").Value = 15.46614604779867
$ws.Range("N9").Value = 16.92277838913891
$ws.Range("O9").Value = 21.98933812791935
$ws.Range("B10").Value = 15.0663111178099
$ws.Range("C10").Value = 9.743100706468852
$ws.Range("D10").Value = 6.455186457539515
$ws.Range("E10").Value = 11.66981086518437
$ws.Range("G10").Value = 3.608759447643167
$ws.Range("M10").Value = 15.95761914476567
$ws.Range("N10").Value = 16.81348375009383
$ws.Range("O10").Value = 22.1066457231871
$ws.Range("B11").Value = 15.45695590965336
$ws.Range("C11").Value = 10.04258245992036
$ws.Range("D11").Value = 6.554820155454527
$ws.Range("E11").Value = 11.78208069487498
$ws.Range("G11").Value = 3.606986236311045
$ws.Range("M11").Value = 16.17992311966194
$ws.Range("N11").Value = 16.76569825909196
$ws.Range("O11").Value = 22.16709155291569
$ws.Range("B12").Value = 15.60257539153104
$ws.Range("C12").Value = 10.15349261315081
$ws.Range("D12").Value = 6.592404929670385
$ws.Range("E12").Value = 11.82493361712786
$ws.Range("G12").Value = 3.606327185753051
$ws.Range("M12").Value = 16.26383707789748
$ws.Range("N12").Value = 16.74787940320751
$ws.Range("O12").Value = 22.1909887508513
$ws.Range("B13").Value = 15.57131788345504
$ws.Range("C13").Value = 10.12971728915953
$ws.Range("D13").Value = 6.584317420003985
$ws.Range("E13").Value = 11.81569004327101
$ws.Range("G13").Value = 3.606468572522036
$ws.Range("M13").Value = 16.24577792623698
$ws.Range("N13").Value = 16.75170474048561
$ws.Range("O13").Value = 22.18579743459864
$ws.Range("B14").Value = 15.46898295383703
$ws.Range("C14").Value = 10.05175718838708
$ws.Range("D14").Value = 6.557915395501037
$ws.Range("E14").Value = 11.78559967004864
$ws.Range("G14").Value = 3.606931767180717
$ws.Range("M14").Value = 16.18683261079794
$ws.Range("N14").Value = 16.76422675882617
$ws.Range("O14").Value = 22.16903746217456
$ws.Range("B15").Value = 15.40599619584235
$ws.Range("C15").Value = 10.00367897953998
$ws.Range("D15").Value = 6.541723425324288
$ws.Range("E15").Value = 11.76721137965514
$ws.Range("G15").Value = 3.607217103488241
$ws.Range("M15").Value = 16.15068957717142
$ws.Range("N15").Value = 16.77193281666193
$ws.Range("O15").Value = 22.15890237566994
$ws.Range("B16").Value = 15.04046392646053
$ws.Range("C16").Value = 9.723179837000652
$ws.Range("D16").Value = 6.448657128291329
$ws.Range("E16").Value = 11.6625233489235
$ws.Range("G16").Value = 3.608877073769808
$ws.Range("M16").Value = 15.94305813909418
$ws.Range("N16").Value = 16.81664541499264
$ws.Range("O16").Value = 22.10283707964756
$ws.Range("B17").Value = 14.81221894765593
$ws.Range("C17").Value = 9.546663359008512
$ws.Range("D17").Value = 6.391349725390548
$ws.Range("E17").Value = 11.59894670639804
$ws.Range("G17").Value = 3.609917617839351
$ws.Range("M17").Value = 15.81529596424625
$ws.Range("N17").Value = 16.84456919812397
$ws.Range("O17").Value = 22.07024966117145
$ws.Range("B18").Value = 14.67950298029996
$ws.Range("C18").Value = 9.44350733576159
$ws.Range("D18").Value = 6.358322497918083
$ws.Range("E18").Value = 11.56262929348476
$ws.Range("G18").Value = 3.61052429409665
$ws.Range("M18").Value = 15.74169589575995
$ws.Range("N18").Value = 16.86081227032134
$ws.Range("O18").Value = 22.05217358945155
$ws.Range("B19").Value = 14.6343246468701
$ws.Range("C19").Value = 9.408301128731674
$ws.Range("D19").Value = 6.347130070008205
$ws.Range("E19").Value = 11.55037701615912
$ws.Range("G19").Value = 3.610731111771417
$ws.Range("M19").Value = 15.71675914629916
$ws.Range("N19").Value = 16.86634320933523
$ws.Range("O19").Value = 22.04616825708336
$ws.Range("B20").Value = 14.83666541714035
$ws.Range("C20").Value = 9.565622450001891
$ws.Range("D20").Value = 6.39745728009203
$ws.Range("E20").Value = 11.60568894149804
$ws.Range("G20").Value = 3.609806003754273
$ws.Range("M20").Value = 15.82890896068013
$ws.Range("N20").Value = 16.84157783317536
$ws.Range("O20").Value = 22.07364965116896
$ws.Range("B21").Value = 15.49910465599224
$ws.Range("C21").Value = 10.0747237828991
$ws.Range("D21").Value = 6.565674540830085
$ws.Range("E21").Value = 11.79442906646991
$ws.Range("G21").Value = 3.606795378988607
$ws.Range("M21").Value = 16.20415420377979
$ws.Range("N21").Value = 16.76054124715147
$ws.Range("O21").Value = 22.17393301855409
$ws.Range("B22").Value = 15.91854792643017
$ws.Range("C22").Value = 10.39289428515854
$ws.Range("D22").Value = 6.674755638245683
$ws.Range("E22").Value = 11.91973765797483
$ws.Range("G22").Value = 3.604900158735111
$ws.Range("M22").Value = 16.44780150984831
$ws.Range("N22").Value = 16.70918995106003
$ws.Range("O22").Value = 22.24534127816229
$ws.Range("B23").Value = 15.69594993645874
$ws.Range("C23").Value = 10.22441466623438
$ws.Range("D23").Value = 6.616628425713471
$ws.Range("E23").Value = 11.85269234199165
$ws.Range("G23").Value = 3.605905071787226
$ws.Range("M23").Value = 16.31793557985463
$ws.Range("N23").Value = 16.73645021451721
$ws.Range("O23").Value = 22.20669648576349
$ws.Range("B24").Value = 14.82561781615506
$ws.Range("C24").Value = 9.55705625959463
$ws.Range("D24").Value = 6.394696302069404
$ws.Range("E24").Value = 11.60264004622379
$ws.Range("G24").Value = 3.609856438143725
$ws.Range("M24").Value = 15.82275497910282
$ws.Range("N24").Value = 16.84292963938369
$ws.Range("O24").Value = 22.07211046246465
$ws.Range("B25").Value = 13.83253207760569
$ws.Range("C25").Value = 8.775015528783504
$ws.Range("D25").Value = 6.152985064068721
$ws.Range("E25").Value = 11.34263795966547
$ws.Range("G25").Value = 3.61443176590233
$ws.Range("M25").Value = 15.28461622542435
$ws.Range("N25").Value = 16.96475166286984
$ws.Range("O25").Value = 21.95247776050605
